# Scheduled market-data refresh: update computed Leve profit columns (H-N)
# across all job sheets to reflect the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 667.25
$ws.Range("J12").Value = 737.5
$ws.Range("L12").Value = 737.5
$ws.Range("N12").Value = -1077.5

$ws.Range("H17").Value = 1676.7727
$ws.Range("J17").Value = 2193.0625
$ws.Range("L17").Value = 6579.1875
$ws.Range("N17").Value = -6915.1875

$ws.Range("H80").Value = 625.0769
$ws.Range("I80").Value = 634.1429
$ws.Range("K80").Value = 1902.4287
$ws.Range("M80").Value = -904.4287000000002

$ws.Range("H83").Value = 625.0769
$ws.Range("I83").Value = 634.1429
$ws.Range("K83").Value = 5707.2861
$ws.Range("M83").Value = -715.2861000000003

$ws.Range("H107").Value = 326.91666
$ws.Range("I107").Value = 215.75
$ws.Range("K107").Value = 215.75
$ws.Range("M107").Value = 1704.25

$ws.Range("H125").Value = 2693.0476
$ws.Range("I125").Value = 2733.4285
$ws.Range("J125").Value = 2612.2856
$ws.Range("K125").Value = 24600.8565
$ws.Range("L125").Value = 23510.5704
$ws.Range("M125").Value = -22140.8565
$ws.Range("N125").Value = -28430.5704

$ws.Range("H132").Value = 11885.6
$ws.Range("I132").Value = 10940.477
$ws.Range("K132").Value = 32821.431
$ws.Range("M132").Value = -30291.431


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1129
$ws.Range("I2").Value = 1129
$ws.Range("K2").Value = 1129
$ws.Range("M2").Value = -1016

$ws.Range("H32").Value = 3032752
$ws.Range("I32").Value = 564.8929
$ws.Range("K32").Value = 564.8929
$ws.Range("M32").Value = -277.8929000000001

$ws.Range("H97").Value = 2598.5
$ws.Range("I97").Value = 2664.6667
$ws.Range("J97").Value = 2400
$ws.Range("K97").Value = 2664.6667
$ws.Range("L97").Value = 2400
$ws.Range("M97").Value = -2168.6667
$ws.Range("N97").Value = -3392

$ws.Range("H116").Value = 1129
$ws.Range("I116").Value = 1129
$ws.Range("K116").Value = 1129
$ws.Range("M116").Value = 1165


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1129
$ws.Range("I3").Value = 1129
$ws.Range("K3").Value = 1129
$ws.Range("M3").Value = -1015

$ws.Range("H86").Value = 3709.8572
$ws.Range("I86").Value = 2239.8333
$ws.Range("J86").Value = 4812.375
$ws.Range("K86").Value = 2239.8333
$ws.Range("L86").Value = 4812.375
$ws.Range("M86").Value = -1116.8333
$ws.Range("N86").Value = -7058.375

$ws.Range("H89").Value = 3709.8572
$ws.Range("I89").Value = 2239.8333
$ws.Range("J89").Value = 4812.375
$ws.Range("K89").Value = 11199.1665
$ws.Range("L89").Value = 24061.875
$ws.Range("M89").Value = -5583.166499999999
$ws.Range("N89").Value = -35293.875

$ws.Range("H99").Value = 1318.091
$ws.Range("I99").Value = 1330.9
$ws.Range("K99").Value = 1330.9
$ws.Range("M99").Value = 167.0999999999999

$ws.Range("H105").Value = 1483.1666
$ws.Range("I105").Value = 1266.4445
$ws.Range("K105").Value = 1266.4445
$ws.Range("M105").Value = 480.5554999999999


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3984.6843
$ws.Range("I31").Value = 3193.7222
$ws.Range("J31").Value = 4349.7437
$ws.Range("K31").Value = 3193.7222
$ws.Range("L31").Value = 4349.7437
$ws.Range("M31").Value = -2898.7222
$ws.Range("N31").Value = -4939.7437

$ws.Range("H34").Value = 3984.6843
$ws.Range("I34").Value = 3193.7222
$ws.Range("J34").Value = 4349.7437
$ws.Range("K34").Value = 3193.7222
$ws.Range("L34").Value = 4349.7437
$ws.Range("M34").Value = -2991.7222
$ws.Range("N34").Value = -4753.7437

$ws.Range("H99").Value = 3560.889
$ws.Range("I99").Value = 3292.5715
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 3292.5715
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -1794.5715
$ws.Range("N99").Value = -7496

$ws.Range("H122").Value = 509.5
$ws.Range("I122").Value = 509.5
$ws.Range("K122").Value = 1528.5
$ws.Range("M122").Value = 921.5

$ws.Range("H126").Value = 3560.889
$ws.Range("I126").Value = 3292.5715
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 9877.7145
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -7407.7145
$ws.Range("N126").Value = -18440

$ws.Range("H132").Value = 2190.7
$ws.Range("I132").Value = 2207.3157
$ws.Range("J132").Value = 1875
$ws.Range("K132").Value = 6621.9471
$ws.Range("L132").Value = 5625
$ws.Range("M132").Value = -4091.9471
$ws.Range("N132").Value = -10685

$ws.Range("H141").Value = 77663
$ws.Range("J141").Value = 77663
$ws.Range("L141").Value = 77663
$ws.Range("N141").Value = -88023


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1762
$ws.Range("I113").Value = 783.3333
$ws.Range("J113").Value = 2181.4285
$ws.Range("K113").Value = 2349.9999
$ws.Range("L113").Value = 6544.2855
$ws.Range("M113").Value = -179.9998999999998
$ws.Range("N113").Value = -10884.2855


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.22222
$ws.Range("I2").Value = 114.6
$ws.Range("J2").Value = 71
$ws.Range("K2").Value = 114.6
$ws.Range("L2").Value = 71
$ws.Range("M2").Value = -1.599999999999994
$ws.Range("N2").Value = -297

$ws.Range("H97").Value = 937.3333
$ws.Range("I97").Value = 900
$ws.Range("J97").Value = 1012
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 1012
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -2004

$ws.Range("H122").Value = 3132
$ws.Range("I122").Value = 1633.3334
$ws.Range("J122").Value = 4256
$ws.Range("K122").Value = 4900.0002
$ws.Range("L122").Value = 12768
$ws.Range("M122").Value = -2450.0002
$ws.Range("N122").Value = -17668


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 996
$ws.Range("J46").Value = 996.6667
$ws.Range("L46").Value = 996.6667
$ws.Range("N46").Value = -1372.6667

$ws.Range("H55").Value = 638.13043
$ws.Range("I55").Value = 664.5
$ws.Range("K55").Value = 664.5
$ws.Range("M55").Value = -491.5


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2127.25
$ws.Range("I126").Value = 1695
$ws.Range("J126").Value = 3769.8
$ws.Range("K126").Value = 5085
$ws.Range("L126").Value = 11309.4
$ws.Range("M126").Value = -2615
$ws.Range("N126").Value = -16249.4

$ws.Range("H132").Value = 1599.0555
$ws.Range("I132").Value = 1448.8572
$ws.Range("K132").Value = 4346.571599999999
$ws.Range("M132").Value = -1816.571599999999

$ws.Range("H136").Value = 2125
$ws.Range("I136").Value = 1735.579
$ws.Range("K136").Value = 5206.737
$ws.Range("M136").Value = -2656.737

